# Timesheet changes by Ruchika(MT2012119)
#
# Fills in the "OFF" (weekend) marker across columns AC:AL and AO for the
# four "leave" rows (28-31), and records the leave-day totals in AM/AN,
# mirroring the formatting already used by the other "OFF" columns (e.g. Q)
# on the same rows. Also nudges the active selection/scroll position the way
# the author's session ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FebruaryMarch 2013")
$ws.Activate()

# Rows 28-31 each get "OFF" stamped (with the same grey style already used
# by column Q on that row) across AC:AL and AO, plus numeric day totals in
# AM/AN.
$offSource = $ws.Cells.Item(28, 17)   # Q28 - already styled/labelled "OFF"

$rows = @(28, 29, 30, 31)
$amanValues = @{
    28 = @(0, 0)
    29 = @(1, 1)
    30 = @(0, 0)
    31 = @(0, 0)
}

foreach ($r in $rows) {
    # AC:AL (columns 29-38) -> "OFF", styled like Q on this row
    $acAl = $ws.Range($ws.Cells.Item($r, 29), $ws.Cells.Item($r, 38))
    $offSource.Copy()
    $acAl.PasteSpecial(-4122)   # xlPasteFormats
    $acAl.Value = "OFF"

    # AO (column 41) -> "OFF", same styling
    $ao = $ws.Cells.Item($r, 41)
    $offSource.Copy()
    $ao.PasteSpecial(-4122)
    $ao.Value = "OFF"

    # AM/AN (columns 39-40) -> numeric totals, formatting unchanged
    $vals = $amanValues[$r]
    $ws.Cells.Item($r, 39).Value = $vals[0]
    $ws.Cells.Item($r, 40).Value = $vals[1]
}

$excel.CutCopyMode = $false

# Restore the view the workbook was left in: scrolled so AG16 is the
# top-left visible cell, with AQ27 selected.
$win = $excel.ActiveWindow
$win.ScrollColumn = 33   # AG
$win.ScrollRow = 16
$ws.Range("AQ27").Select()
